# create matriz product and espTec
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("PRODUCTO")

# Duplicate PRODUCTO (keeps workbook/sheet defaults clean, no stray
# baseColWidth like a brand-new Worksheets.Add() sheet would have) and turn
# the copy into the new ESP_TECNICAS sheet, positioned right after PRODUCTO.
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item("PRODUCTO (2)")
$ws2.Name = "ESP_TECNICAS"
$ws2.Cells.Clear()

# --- PRODUCTO: rewrite as the new 2-row product "matriz" (was 8 spec rows
# plus 2 padding rows in A1:G10; now just A2:E3). ---
$ws1.Cells.Clear()

function Set-TextCell($range, $value) {
    # Force values that look numeric (leading zeros, etc.) to be stored as
    # text/shared-strings instead of being coerced to a number, without
    # leaving a lingering cell style behind.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextCell $ws1.Range("A2") "01000101"
$ws1.Range("B2").Value = 1
$ws1.Range("C2").Value = "MEGIMPERU"
$ws1.Range("D2").Value = "Gancho Organizador para Reposacabezas de Asiento"
$ws1.Range("E2").Value = 1

Set-TextCell $ws1.Range("A3") "01000101"
$ws1.Range("B3").Value = 2
$ws1.Range("C3").Value = "TOP GAN"
$ws1.Range("D3").Value = "Piso para Auto TOP GAM Negro"
$ws1.Range("E3").Value = 1

# --- ESP_TECNICAS: the technical-specification rows that used to live on
# PRODUCTO (rows 1-8), renumbered to rows 2-9. ---
$specs = @(
    @(1, "Tipo", "Organizador"),
    @(2, "Color", "Negro"),
    @(3, "Material", "ABS"),
    @(4, "Características", "Organiza bolsas de mano, bolsos de mano, mochilas, paraguas y más"),
    @(5, "Modelo", "Gancho reposacabezas"),
    @(6, "Alto (cm)", "6"),
    @(7, "Ancho (cm)", "11"),
    @(8, "Incluye", "1 gancho de doble percha")
)

$r = 2
foreach ($spec in $specs) {
    Set-TextCell $ws2.Cells.Item($r, 1) "01000101"
    $ws2.Cells.Item($r, 2).Value = 1
    $ws2.Cells.Item($r, 3).Value = $spec[0]
    $ws2.Cells.Item($r, 4).Value = $spec[1]
    Set-TextCell $ws2.Cells.Item($r, 5) $spec[2]
    $ws2.Cells.Item($r, 6).Value = 0
    $ws2.Cells.Item($r, 7).Value = 1
    $r++
}

# Re-select PRODUCTO as the active sheet/tab (fresh lookups so the
# activation actually sticks on the final workbook state).
$final = $wb.Worksheets.Item("PRODUCTO")
$final.Activate()
$null = $final.Range("A1").Select()
